$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The stimulus set was regenerated with updated distance/size codes:
#   D80 -> D86, D51 -> D55, D64 -> D69, S30 -> S31
# These codes appear as substrings inside many string values across several
# columns (Condition, Filename_Left, Filename_Right, Distance, Size), so do a
# whole-sheet text substitution for each old/new pair.

$ws.Cells.Replace("D80", "D86")
$ws.Cells.Replace("D51", "D55")
$ws.Cells.Replace("D64", "D69")
$ws.Cells.Replace("S30", "S31")

Write-Output "renamed distance/size codes"
